# Add Bernard lake station
# - Rename "Bernard" sheet to "Bernard_lake"
# - Update coordinates and a couple of cell values on that sheet
# - Add a new shared string "Thermistors chain 1" and use it on the
#   Thermistors sheet (B2), replacing the old reference to "Thermistors chain 2"
# - Update selections on several sheets and move the active tab from
#   "Berge" to "Thermistors"

$wb = $excel.ActiveWorkbook

# Rename the "Bernard" worksheet to "Bernard_lake"
$wsBernard = $wb.Worksheets.Item("Bernard")
$wsBernard.Name = "Bernard_lake"

# Update values on the Bernard_lake sheet
$wsBernard.Range("G2").Value = 50.8614547
$wsBernard.Range("G3").Value = -63.3898409
$wsBernard.Range("C17").Value = -37.5
$wsBernard.Range("D17").Value = 5.5

# Update the Thermistors sheet label text (adds a new shared string)
$wsThermistors = $wb.Worksheets.Item("Thermistors")
$wsThermistors.Range("B2").Value = "Thermistors chain 1"

# Update selections (and, implicitly, which sheet/tab is active).
# Selecting a range on a sheet makes that sheet active, so we select the
# sheets that keep a (non-default) selection first and finish with the
# sheet that should end up as the active tab.
$wsBerge = $wb.Worksheets.Item("Berge")
[void]$wsBerge.Range("C17").Select()

[void]$wsBernard.Range("C5").Select()

$wsForet = $wb.Worksheets.Item("Foret")
[void]$wsForet.Range("F27").Select()

[void]$wsThermistors.Range("B3").Select()
